# AutoCommit_12 апреля 2024 г. 10:04:43_SibNout2023
#
# Sibirev I. V. gradebook: row 26 (Сушко Артур) gets marked for ДЗ_1..ДЗ_3
# (columns C:E) — was 0/0/0 (still highlighted "not yet graded" green),
# becomes 5/5/5 with the "graded" (no-fill) look used elsewhere in the
# sheet. The row's sum (column J) recalculates automatically from the
# existing shared formula. Finally, the active selection moves to F29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the "graded" cell style (no fill) already used by e.g. C4 on this
# sheet, instead of the "ungraded" green-fill style C29:E29 currently have.
# Copy format only so the donor cell's value (4) is left untouched.
$formatDonor = $ws.Range("C4")
$target = $ws.Range("C29:E29")
$formatDonor.Copy()
$target.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the grades.
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 5

# Leave the selection on F29, as in the saved workbook.
$ws.Range("F29").Select()
